# The sheet that used to act as the chart "template" carried a couple of
# leftover, manually-pasted helper percentages in U5:U6 (249/257 and
# 32/249). Clear just the contents (formula + cached value) so the
# existing number formatting / style on those cells is preserved, then
# leave that range selected - matching how the workbook was left after
# the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U5:U6").ClearContents()
$ws.Range("U5:U6").Select()
